$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 14000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 14000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 14000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -14972
$ws.Range("H98").Value = 703.4667
$ws.Range("I98").Value = 700.13794
$ws.Range("K98").Value = 700.13794
$ws.Range("M98").Value = 797.86206
$ws.Range("H109").Value = 23299.7
$ws.Range("J109").Value = 23299.7
$ws.Range("L109").Value = 23299.7
$ws.Range("N109").Value = -26073.7
$ws.Range("H122").Value = 703.4667
$ws.Range("I122").Value = 700.13794
$ws.Range("K122").Value = 2100.41382
$ws.Range("M122").Value = 349.5861800000002
$ws.Range("H129").Value = 1000328.94
$ws.Range("J129").Value = 1158633.8
$ws.Range("L129").Value = 3475901.4
$ws.Range("N129").Value = -3485901.4
$ws.Range("H140").Value = 69000
$ws.Range("J140").Value = 69000
$ws.Range("L140").Value = 69000
$ws.Range("N140").Value = -79360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 7406.6665
$ws.Range("J44").Value = 7406.6665
$ws.Range("L44").Value = 7406.6665
$ws.Range("N44").Value = -8382.6665
$ws.Range("H55").Value = 10908.75
$ws.Range("I55").Value = 11000
$ws.Range("J55").Value = 10895.714
$ws.Range("K55").Value = 11000
$ws.Range("L55").Value = 10895.714
$ws.Range("M55").Value = -10685
$ws.Range("N55").Value = -11525.714
$ws.Range("H80").Value = 21463.334
$ws.Range("J80").Value = 21463.334
$ws.Range("L80").Value = 21463.334
$ws.Range("N80").Value = -23459.334
$ws.Range("H83").Value = 21463.334
$ws.Range("J83").Value = 21463.334
$ws.Range("L83").Value = 64390.00199999999
$ws.Range("N83").Value = -74374.002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24180.637
$ws.Range("I82").Value = 5226
$ws.Range("J82").Value = 31288.625
$ws.Range("K82").Value = 5226
$ws.Range("L82").Value = 31288.625
$ws.Range("M82").Value = -4843
$ws.Range("N82").Value = -32054.625
$ws.Range("H85").Value = 24180.637
$ws.Range("I85").Value = 5226
$ws.Range("J85").Value = 31288.625
$ws.Range("K85").Value = 5226
$ws.Range("L85").Value = 31288.625
$ws.Range("M85").Value = -3900
$ws.Range("N85").Value = -33940.625
$ws.Range("H99").Value = 1884.3846
$ws.Range("I99").Value = 1926.6666
$ws.Range("J99").Value = 1871.7
$ws.Range("K99").Value = 1926.6666
$ws.Range("L99").Value = 1871.7
$ws.Range("M99").Value = -428.6666
$ws.Range("N99").Value = -4867.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 404
$ws.Range("J34").Value = 1500
$ws.Range("L34").Value = 4500
$ws.Range("N34").Value = -4668
$ws.Range("H62").Value = 3672.4
$ws.Range("I62").Value = 776.3333
$ws.Range("J62").Value = 4913.5713
$ws.Range("K62").Value = 2328.9999
$ws.Range("L62").Value = 14740.7139
$ws.Range("M62").Value = -1642.9999
$ws.Range("N62").Value = -16112.7139
$ws.Range("H65").Value = 3672.4
$ws.Range("I65").Value = 776.3333
$ws.Range("J65").Value = 4913.5713
$ws.Range("K65").Value = 6986.9997
$ws.Range("L65").Value = 44222.14169999999
$ws.Range("M65").Value = -3554.9997
$ws.Range("N65").Value = -51086.14169999999
$ws.Range("H69").Value = 2333.3333
$ws.Range("J69").Value = 2333.3333
$ws.Range("L69").Value = 6999.999899999999
$ws.Range("N69").Value = -8621.999899999999
$ws.Range("H72").Value = 2333.3333
$ws.Range("J72").Value = 2333.3333
$ws.Range("L72").Value = 20999.9997
$ws.Range("N72").Value = -29111.9997
$ws.Range("H86").Value = 1466.6666
$ws.Range("I86").Value = 400
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -14
$ws.Range("N86").Value = -8372
$ws.Range("H89").Value = 1466.6666
$ws.Range("I89").Value = 400
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 3600
$ws.Range("L89").Value = 18000
$ws.Range("M89").Value = 2328
$ws.Range("N89").Value = -29856
$ws.Range("H131").Value = 671.88
$ws.Range("I131").Value = 293.95834
$ws.Range("J131").Value = 791.2237
$ws.Range("K131").Value = 881.8750200000001
$ws.Range("L131").Value = 2373.6711
$ws.Range("M131").Value = 4158.12498
$ws.Range("N131").Value = -12453.6711

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 123824.414
$ws.Range("I70").Value = 226822.89
$ws.Range("K70").Value = 226822.89
$ws.Range("M70").Value = -226552.89
$ws.Range("H73").Value = 123824.414
$ws.Range("I73").Value = 226822.89
$ws.Range("K73").Value = 226822.89
$ws.Range("M73").Value = -225886.89
$ws.Range("H95").Value = 16896
$ws.Range("J95").Value = 16896
$ws.Range("L95").Value = 16896
$ws.Range("N95").Value = -22388
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 2556.074
$ws.Range("I132").Value = 2185.1
$ws.Range("J132").Value = 3616
$ws.Range("K132").Value = 6555.299999999999
$ws.Range("L132").Value = 10848
$ws.Range("M132").Value = -4025.299999999999
$ws.Range("N132").Value = -15908

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1821.6428
$ws.Range("I22").Value = 5899.5
$ws.Range("J22").Value = 1142
$ws.Range("K22").Value = 5899.5
$ws.Range("L22").Value = 1142
$ws.Range("M22").Value = -5604.5
$ws.Range("N22").Value = -1732
$ws.Range("H27").Value = 1821.6428
$ws.Range("I27").Value = 5899.5
$ws.Range("J27").Value = 1142
$ws.Range("K27").Value = 5899.5
$ws.Range("L27").Value = 1142
$ws.Range("M27").Value = -5792.5
$ws.Range("N27").Value = -1356
$ws.Range("H101").Value = 29000
$ws.Range("J101").Value = 29000
$ws.Range("L101").Value = 29000
$ws.Range("N101").Value = -35490
$ws.Range("H122").Value = 4999.6665
$ws.Range("I122").Value = 4999.6665
$ws.Range("K122").Value = 14998.9995
$ws.Range("M122").Value = -12548.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 36307.5
$ws.Range("J105").Value = 36307.5
$ws.Range("L105").Value = 36307.5
$ws.Range("N105").Value = -43295.5
$ws.Range("H132").Value = 1134.7693
$ws.Range("I132").Value = 704.36365
$ws.Range("K132").Value = 2113.09095
$ws.Range("M132").Value = 416.9090500000002
